# "switch levels 3 and 4"
#
# - levels sheet: row 4 (level 3) and row 5 (level 4) swap their name/days
#   values, and get highlighted green.
# - techniques sheet: the "level" column (D) is re-numbered for most rows.
# - contents sheet: one "standards" count (D72) changes.
# - actions sheet: a batch of "standards" counts (F59:F71) change, and one
#   action-type letter (G70) changes.
# - level_techniquesList sheet: level references in column B are swapped
#   between 3 and 4 (old level-3 rows -> 4, highlighted yellow; old
#   level-4 rows -> 3, highlighted green).

$wb = $excel.ActiveWorkbook

$yellow = 65535      # RGB(255,255,0)
$green  = 5296274    # RGB(146,208,80)

# ---------------------------------------------------------------------
# level_techniquesList: renumber level ids, yellow first so it becomes
# the first newly-created fill/style (matches fillId=2/style 28), then
# green becomes the second (fillId=3/style 29).
# ---------------------------------------------------------------------
$wsLT = $wb.Worksheets.Item("level_techniquesList")

$wsLT.Range("B8:B12").Value2 = 4
$wsLT.Range("B8:B12").Interior.Color = $yellow

$wsLT.Range("B13:B17").Value2 = 3
$wsLT.Range("B13:B17").Interior.Color = $green

$wsLT.Activate()
$wsLT.Range("B13:B17").Select()

# ---------------------------------------------------------------------
# levels: swap the level-3 / level-4 rows' name + day-count, highlight
# both rows green.
# ---------------------------------------------------------------------
$wsLevels = $wb.Worksheets.Item("levels")

$name3 = $wsLevels.Range("B4").Value2
$name4 = $wsLevels.Range("B5").Value2
$days3 = $wsLevels.Range("C4").Value2
$days4 = $wsLevels.Range("C5").Value2

$wsLevels.Range("B4").Value2 = $name4
$wsLevels.Range("C4").Value2 = $days4
$wsLevels.Range("B5").Value2 = $name3
$wsLevels.Range("C5").Value2 = $days3

$wsLevels.Range("B4:B5").Interior.Color = $green

$wsLevels.Activate()
$wsLevels.Range("C6").Select()

# ---------------------------------------------------------------------
# techniques: renumber the level column (D).
# ---------------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("techniques")

$wsTech.Range("D2:D5").Value2 = 3
$wsTech.Range("D6:D9").Value2 = 1
$wsTech.Range("D10").Value2 = 3
$wsTech.Range("D12").Value2 = 2
$wsTech.Range("D14:D17").Value2 = 3
$wsTech.Range("D18").Value2 = 0

$wsTech.Columns.Item(2).ColumnWidth = 59.3
$wsTech.Columns.Item(3).ColumnWidth = 9

$wsTech.Activate()
$wsTech.Range("D20").Select()

# ---------------------------------------------------------------------
# contents: bump the standards count on row 72.
# ---------------------------------------------------------------------
$wsContents = $wb.Worksheets.Item("contents")

$wsContents.Range("D72").Value2 = 11

$wsContents.Activate()
$wsContents.Range("D1").Select()

# ---------------------------------------------------------------------
# actions: bump the standards counts for rows 59:71, fix the action-type
# letter on row 70.
# ---------------------------------------------------------------------
$wsActions = $wb.Worksheets.Item("actions")

$wsActions.Range("F59:F71").Value2 = 15
$wsActions.Range("G70").Value2 = "T"

$wsActions.Activate()
$wsActions.Range("C58").Select()

# ---------------------------------------------------------------------
# Final active sheet is "levels" (matches the workbook's active tab).
# ---------------------------------------------------------------------
$wsLevels.Activate()
$wsLevels.Range("C6").Select()
